# "Doing Updates for Financials"
# Update the latest fiscal-year (column D) figures on the JKHY sheet with
# refreshed financial data, plus a handful of corrected prior-year (E:J)
# Capital Expenditures figures on row 91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Income Statement - Non Recurring
$ws.Range("D14").Value = -1600

# Balance Sheet - Net Receivables
$ws.Range("D43").Value = 318900

# Balance Sheet - Other Current Assets
$ws.Range("D45").Value = 123200

# Balance Sheet - Total Current Assets
$ws.Range("D46").Value = 473600

# Balance Sheet - Property Plant and Equipment
$ws.Range("D48").Value = 937900

# Balance Sheet - Goodwill
$ws.Range("D49").Value = 1533300

# Balance Sheet - Other Assets
$ws.Range("D52").Value = 185200

# Balance Sheet - Total Assets
$ws.Range("D54").Value = 2037200

# Balance Sheet - Other Current Liabilities
$ws.Range("D59").Value = 441200

# Balance Sheet - Total Current Liabilities
$ws.Range("D60").Value = 475700

# Balance Sheet - Other Liabilities
$ws.Range("D62").Value = 428300

# Balance Sheet - Total Liabilities
$ws.Range("D66").Value = 714400

# Balance Sheet - Retained Earnings
$ws.Range("D72").Value = 1912900

# Balance Sheet - Total Stockholder Equity
$ws.Range("D76").Value = 1322800

# Cash Flow Statement - Capital Expenditures (all years revised)
$ws.Range("D91").Value = -40100
$ws.Range("E91").Value = -41900
$ws.Range("F91").Value = -56300
$ws.Range("G91").Value = -54400
$ws.Range("H91").Value = -33200
$ws.Range("I91").Value = -46300
$ws.Range("J91").Value = -41400
